# Add new template rows 047, 048, 049 (plus their related 046 row) to the
# "Templates" sheet, as described by the commit "new templates : 047, 048, 049".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Templates")

# Row 47 -> Template 046
$ws.Range("A47").Value = "046"
$ws.Range("B47").Value = 386
$ws.Range("C47").Value = 466
$ws.Range("D47").Value = 684
$ws.Range("E47").Value = 524
$ws.Range("F47").Value = "046"
$ws.Range("G47").Value = "unable to play dialog - max number of player -> join"

# Row 48 -> Template 047
$ws.Range("A48").Value = "047"
$ws.Range("B48").Value = 886
$ws.Range("C48").Value = 681
$ws.Range("D48").Value = 960
$ws.Range("E48").Value = 723
$ws.Range("F48").Value = "047"
$ws.Range("G48").Value = "ok button - unable to play dialog"

# Row 49 -> Template 048
$ws.Range("A49").Value = "048"
$ws.Range("B49").Value = 495
$ws.Range("C49").Value = 152
$ws.Range("D49").Value = 1342
$ws.Range("E49").Value = 196
$ws.Range("F49").Value = "048"
$ws.Range("G49").Value = "match condition have not met dialog - join"

# Row 50 -> Template 049
$ws.Range("A50").Value = "049"
$ws.Range("B50").Value = 463
$ws.Range("C50").Value = 826
$ws.Range("D50").Value = 615
$ws.Range("E50").Value = 893
$ws.Range("F50").Value = "020"
$ws.Range("G50").Value = "close button - match condition have not met dialog"

# Match the author's final selection/active cell after entering the new data.
$ws.Range("E50").Select()
